$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6041439175605774
$ws.Range("B1").Value = 0.7307738661766052
$ws.Range("C1").Value = 0.9839804768562317
$ws.Range("D1").Value = 2.917953014373779
$ws.Range("E1").Value = 4.959304332733154
